$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") '26.133.23'
$ws.Range("E2").Value = '  -0.72%  '

Set-TextCell $ws.Range("D3") '1.668.73'
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("E4").Value = '  -0.55%  '

Set-TextCell $ws.Range("D5") '210.56'
$ws.Range("E5").Value = '  -3.39%  '

$ws.Range("E6").Value = '  -2.85%  '

$ws.Range("E7").Value = '  -0.54%  '

Set-TextCell $ws.Range("D8") '0.2633'
$ws.Range("E8").Value = '  -3.85%  '

Set-TextCell $ws.Range("D9") '0.06300'
$ws.Range("E9").Value = '  -2.42%  '

Set-TextCell $ws.Range("D10") '21.20'
$ws.Range("E10").Value = '  -2.32%  '

Set-TextCell $ws.Range("D11") '0.07568'
$ws.Range("E11").Value = '  -1.18%  '

Set-TextCell $ws.Range("D12") '1.669.40'
$ws.Range("E12").Value = '  -1.98%  '

Set-TextCell $ws.Range("D13") '4.445'
$ws.Range("E13").Value = '  -2.12%  '

Set-TextCell $ws.Range("D14") '0.5572'
$ws.Range("E14").Value = '  -3.91%  '

Set-TextCell $ws.Range("D15") '66.90'
$ws.Range("E15").Value = '  -0.02%  '

Set-TextCell $ws.Range("D16") '0.000007926'
$ws.Range("E16").Value = '  -5.82%  '

Set-TextCell $ws.Range("D17") '26.162.05'
$ws.Range("E17").Value = '  -0.79%  '

$ws.Range("E18").Value = '  -0.57%  '

Set-TextCell $ws.Range("D19") '4.746'
$ws.Range("E19").Value = '  -3.55%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws.Range("D20") '10.38'
$ws.Range("E20").Value = '  -4.53%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell $ws.Range("D21") '186.71'
$ws.Range("E21").Value = '  -1.92%  '

Set-TextCell $ws.Range("D22") '6.183'
$ws.Range("E22").Value = '  -1.39%  '

Set-TextCell $ws.Range("D23") '1.005'
$ws.Range("E23").Value = '  -0.49%  '

Set-TextCell $ws.Range("D24") '149.57'
$ws.Range("E24").Value = '  +0.24%  '

Set-TextCell $ws.Range("D25") '0.1252'
$ws.Range("E25").Value = '  -2.50%  '

Set-TextCell $ws.Range("D26") '7.515'
$ws.Range("E26").Value = '  -4.37%  '

$ws.Range("E27").Value = '  +0.94%  '

Set-TextCell $ws.Range("D28") '0.06260'
$ws.Range("E28").Value = '  -0.77%  '

$ws.Range("E29").Value = '  -2.02%  '

Set-TextCell $ws.Range("D30") '1.281'
$ws.Range("E30").Value = '  -3.28%  '

Set-TextCell $ws.Range("D31") '3.520'
$ws.Range("E31").Value = '  -2.35%  '

Set-TextCell $ws.Range("D32") '3.421'
$ws.Range("E32").Value = '  -4.70%  '

Set-TextCell $ws.Range("D33") '1.632'
$ws.Range("E33").Value = '  -2.97%  '

$ws.Range("E34").Value = '  -3.19%  '

Set-TextCell $ws.Range("D35") '0.6062'
$ws.Range("E35").Value = '  -2.08%  '

Set-TextCell $ws.Range("D36") '2.412'
$ws.Range("E36").Value = '  -0.19%  '

Set-TextCell $ws.Range("D37") '2.731'
$ws.Range("E37").Value = '  -1.36%  '

Set-TextCell $ws.Range("D38") '6.139'
$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range("D39") '0.01617'
$ws.Range("E39").Value = '  -2.46%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws.Range("D40") '1.101.41'
$ws.Range("E40").Value = '  -1.05%  '

Set-TextCell $ws.Range("D41") '0.8738'
$ws.Range("E41").Value = '  -1.32%  '

$ws.Range("E42").Value = '  -0.99%  '

$ws.Range("E43").Value = '  -0.72%  '

Set-TextCell $ws.Range("D44") '1.823.50'
$ws.Range("E44").Value = '  -1.11%  '

$ws.Range("E45").Value = '  -1.44%  '

Set-TextCell $ws.Range("D46") '55.33'
$ws.Range("E46").Value = '  -4.13%  '

Set-TextCell $ws.Range("D47") '1.005'
$ws.Range("E47").Value = '  +0.07%  '

Set-TextCell $ws.Range("D48") '8.038'
$ws.Range("E48").Value = '  -1.93%  '

Set-TextCell $ws.Range("D49") '0.05236'
$ws.Range("E49").Value = '  -0.86%  '

Set-TextCell $ws.Range("D50") '0.4247'
$ws.Range("E50").Value = '  -1.22%  '

Set-TextCell $ws.Range("D51") '5.990'
